# Adds a new "2022-Q4" data sheet ahead of the existing quarters, shifting
# the older quarterly tabs (and the totals table) down to make room.
#
# Resulting tab order:
#   总计      - totals table, gets a new top row for 2022-Q4
#   2022-Q4   - brand-new quarter numbers (reuses the sheet that used to be "2022-Q3")
#   2022-Q3   - the former "2022-Q3" numbers, kept as history (new duplicate tab)
#   2022-Q2   - unchanged, just shifted one tab to the right

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item("总计")
$q3Sheet    = $wb.Worksheets.Item("2022-Q3")

# --- 1) Duplicate "2022-Q3" (values + styles) right after itself; this copy
#        keeps the historical "2022-Q3" figures. ------------------------------
$q3Sheet.Copy($null, $q3Sheet)
$q3Copy = $wb.Worksheets.Item($q3Sheet.Index + 1)

# Free up the "2022-Q3" name on the original sheet before reusing it below.
$q3Sheet.Name = "2022-Q4"
$q3Copy.Name = "2022-Q3"

# --- 2) Overwrite the renamed sheet's fund-holding figures with the new
#        quarter's numbers (same fund code/name rows, updated stats). --------
$q3Sheet.Range("D2:G4").NumberFormat = "@"

$q3Sheet.Range("D2").Value = "12.08"
$q3Sheet.Range("E2").Value = "89.83"
$q3Sheet.Range("F2").Value = "4.81"
$q3Sheet.Range("G2").Value = "0.5810"
$q3Sheet.Range("H2").Value = 7

$q3Sheet.Range("D3").Value = "12.08"
$q3Sheet.Range("E3").Value = "89.83"
$q3Sheet.Range("F3").Value = "4.81"
$q3Sheet.Range("G3").Value = "0.5810"
$q3Sheet.Range("H3").Value = 7

$q3Sheet.Range("D4").Value = "12.08"
$q3Sheet.Range("E4").Value = "89.83"
$q3Sheet.Range("F4").Value = "4.81"
$q3Sheet.Range("G4").Value = "0.5810"
$q3Sheet.Range("H4").Value = 7

# --- 3) Push the totals sheet's existing rows down one and insert the new
#        2022-Q4 total as the new row 2, preserving each row's formatting. ---
$totalSheet.Range("B3:D3").Copy()
$totalSheet.Range("B4:D4").PasteSpecial(-4163)   # xlPasteAll
$totalSheet.Range("B2:D2").Copy()
$totalSheet.Range("B3:D3").PasteSpecial(-4163)   # xlPasteAll

$totalSheet.Range("B2").Value = "2022-Q4"
$totalSheet.Range("C2").Value = 3
$totalSheet.Range("D2").Value = 1.74

# Column A is just the 0-based row index; restyle to match and renumber.
$totalSheet.Range("A2").Copy()
$totalSheet.Range("A3").PasteSpecial(-4122)      # xlPasteFormats
$totalSheet.Range("A4").PasteSpecial(-4122)      # xlPasteFormats
$totalSheet.Range("A2").Value = 0
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2

# Keep the workbook's active tab where it started (总计 / first sheet).
$totalSheet.Activate()
